$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.405.42"
$ws.Range("E2").Value = "'  -8.05%  "
$ws.Range("D3").Value = "'1.680.56"
$ws.Range("E3").Value = "'  -6.81%  "
$ws.Range("D4").Value = "'1.007"
$ws.Range("E4").Value = "'  +0.44%  "
$ws.Range("D5").Value = "'216.82"
$ws.Range("E5").Value = "'  -6.41%  "
$ws.Range("D6").Value = "'1.007"
$ws.Range("E6").Value = "'  +0.35%  "
$ws.Range("D7").Value = "'0.4965"
$ws.Range("E7").Value = "'  -16.57%  "
$ws.Range("D8").Value = "'0.2603"
$ws.Range("E8").Value = "'  -6.32%  "
$ws.Range("D9").Value = "'21.66"
$ws.Range("E9").Value = "'  -7.18%  "
$ws.Range("D10").Value = "'0.06137"
$ws.Range("E10").Value = "'  -10.20%  "
$ws.Range("D11").Value = "'0.07296"
$ws.Range("E11").Value = "'  -3.25%  "
$ws.Range("D12").Value = "'1.711.28"
$ws.Range("E12").Value = "'  -5.04%  "
$ws.Range("D13").Value = "'4.419"
$ws.Range("E13").Value = "'  -7.98%  "
$ws.Range("D14").Value = "'0.5714"
$ws.Range("E14").Value = "'  -8.61%  "
$ws.Range("D15").Value = "'1.913.16"
$ws.Range("E15").Value = "'  -6.62%  "
$ws.Range("D16").Value = "'0.000008190"
$ws.Range("E16").Value = "'  -11.66%  "
$ws.Range("D17").Value = "'64.24"
$ws.Range("E17").Value = "'  -14.76%  "
$ws.Range("D18").Value = "'26.468.66"
$ws.Range("E18").Value = "'  -7.78%  "
$ws.Range("D19").Value = "'4.983"
$ws.Range("E19").Value = "'  -8.90%  "
$ws.Range("D20").Value = "'1.006"
$ws.Range("E20").Value = "'  +0.28%  "
$ws.Range("D21").Value = "'10.69"
$ws.Range("E21").Value = "'  -6.66%  "
$ws.Range("D22").Value = "'182.83"
$ws.Range("E22").Value = "'  -13.23%  "
$ws.Range("D23").Value = "'6.154"
$ws.Range("E23").Value = "'  -10.26%  "
$ws.Range("D24").Value = "'1.009"
$ws.Range("E24").Value = "'  +0.51%  "
$ws.Range("D25").Value = "'144.39"
$ws.Range("E25").Value = "'  -6.40%  "
$ws.Range("D26").Value = "'7.512"
$ws.Range("E26").Value = "'  -4.22%  "
$ws.Range("D27").Value = "'0.1129"
$ws.Range("E27").Value = "'  -11.53%  "
$ws.Range("D28").Value = "'15.39"
$ws.Range("E28").Value = "'  -6.26%  "
$ws.Range("D29").Value = "'1.311"
$ws.Range("E29").Value = "'  -8.49%  "
$ws.Range("D30").Value = "'0.05622"
$ws.Range("E30").Value = "'  -9.25%  "
$ws.Range("E31").Value = "'  -6.79%  "
$ws.Range("D32").Value = "'3.468"
$ws.Range("E32").Value = "'  -8.29%  "
$ws.Range("D33").Value = "'3.450"
$ws.Range("E33").Value = "'  -7.95%  "
$ws.Range("D34").Value = "'1.626"
$ws.Range("E34").Value = "'  -5.37%  "
$ws.Range("D35").Value = "'1.001"
$ws.Range("E35").Value = "'  -5.68%  "
$ws.Range("D36").Value = "'2.372"
$ws.Range("E36").Value = "'  -4.97%  "
$ws.Range("D37").Value = "'0.5870"
$ws.Range("E37").Value = "'  -8.31%  "
$ws.Range("D38").Value = "'2.632"
$ws.Range("E38").Value = "'  -3.31%  "
$ws.Range("D39").Value = "'0.01582"
$ws.Range("E39").Value = "'  -7.65%  "
$ws.Range("B40").Value = "'Maker"
$ws.Range("C40").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "'1.069.57"
$ws.Range("E40").Value = "'  -5.68%  "
$ws.Range("B41").Value = "'FraxShare"
$ws.Range("C41").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'5.903"
$ws.Range("E41").Value = "'  -8.17%  "
$ws.Range("D42").Value = "'0.8507"
$ws.Range("E42").Value = "'  -2.10%  "
$ws.Range("E43").Value = "'  -0.08%  "
$ws.Range("D44").Value = "'98.15"
$ws.Range("E44").Value = "'  -2.38%  "
$ws.Range("D45").Value = "'1.841.80"
$ws.Range("E45").Value = "'  -6.15%  "
$ws.Range("D46").Value = "'56.26"
$ws.Range("E46").Value = "'  -7.15%  "
$ws.Range("E47").Value = "'  -6.06%  "
$ws.Range("D48").Value = "'1.004"
$ws.Range("E48").Value = "'  +0.05%  "
$ws.Range("D49").Value = "'8.076"
$ws.Range("E49").Value = "'  -3.31%  "
$ws.Range("E50").Value = "'  -3.54%  "
$ws.Range("D51").Value = "'0.05184"
$ws.Range("E51").Value = "'  -5.30%  "
